# Update "想去人数" (number of people interested) values that changed
# between publishing runs of the 丽水-漫展信息 scraper output.
#
# Sheet "展览" (exhibitions) and sheet "全部类型" (all types) both contain
# the same two event rows; row 2's F value goes 68 -> 69 and row 3's F
# value goes 2 -> 4.

$wb = $excel.ActiveWorkbook

foreach ($sheetName in @("展览", "全部类型")) {
    $ws = $wb.Worksheets.Item($sheetName)
    $ws.Range("F2").Value = 69
    $ws.Range("F3").Value = 4
}
